$wb = $excel.ActiveWorkbook

# 1. Update the "psdquery" row's query text in the queries sheet: it used to
#    reference the old "projectStructureDetails" table, now it should query
#    "projectLevelDetails" instead.
$wsQueries = $wb.Worksheets.Item("queries")
$wsQueries.Range("B2").Value = "select * from projectLevelDetails"

# 2. The projectLevelDetails sheet keeps its own remembered cell selection,
#    which moved to G10.
$wsProjectLevelDetails = $wb.Worksheets.Item("projectLevelDetails")
$null = $wsProjectLevelDetails.Range("G10").Select()

# 3. The "queries" sheet becomes the active/selected tab (previously
#    "testCasesTestNG" was active).
$wsQueries.Activate()

# 4. ... and its remembered selection moves to B8.
$null = $wsQueries.Range("B8").Select()
